$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the admin-log row as literal text values (matches inlineStr in the
# target XML -- dates/booleans must stay as plain strings, not get coerced).
$ws.Range("A1").Value = "2025-07-23 09:01:35"
$ws.Range("B1").Value = "add-user"
$ws.Range("C1").Value = "new-organization97"
$ws.Range("D1").Value = "firstteam"
$ws.Range("F1").Value = "Vignesh2122"
$ws.Range("G1").Value = "pull"
$ws.Range("I1").Value = "'False"
